$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$cell = $t.Cell(11, 2)
$cell.Range.Delete()
$p = $cell.Range.Paragraphs.Item(1)
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>29/08 – BLC – PERIPHERAL–RIGHT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – UC – CATHETER SPECIMEN URINE (CSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – UC – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic usage advice, no microbiological findings.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – SARS–CoV–2 RNA – Negative</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – RESP. CULT AND MICRO – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Candida albicans, Neisseria spp. &lt;10^4 ORGs/mL.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – RESP. CULT AND MICRO – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Candida and Neisseria spp. below threshold.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>28/08 – SARS CORONAVIRUS–2 PCR – **No clear Result** +</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: MS2 control detected; no pathogen result.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>26/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>26/08 – UC – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic advice; no microbiological findings stated.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>26/08 – UC – CATHETER SPECIMEN URINE (CSU)PERIPHERAL–LEFT NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>26/08 – RESP. CULT AND MICRO – **No clear Result** +</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Coag–neg Staph, Strep viridans &lt;10^4.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Reasoning:** Both organisms are reported below the diagnostic threshold (10^4 CFU/mL for BAL cultures), making clinical significance uncertain. While they are potential pathogens, the low colony count does not definitively indicate infection (commonly considered contamination or colonization at this level).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>26/08 – RESP. CULT AND MICRO – **No clear Result** +</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Coag–neg Staph, Strep viridans &lt;10^4.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>26/08 – C difficile(GDH) – Negative</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – UC – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic recommendations, no organism data.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – UC – CATHETER SPECIMEN URINE (CSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – RESP. CULT AND MICRO – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Neisseria and Strep. viridans &lt;10^4.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>*Rationale: While organisms are detected, their quantities (&lt;10^4 CFU/mL) fall below typical thresholds for clinical significance in BAL cultures, making the result non–definitive.*</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – RESP. CULT AND MICRO – **No clear Result**</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Neisseria spp., Strep viridans &lt;10^4.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Explanation:** The reported organisms are below the diagnostic threshold (10^4 CFU/mL), making the result inconclusive for infection.</w:t></w:r></w:p>'
$p.Range.InsertXML($xmlFrag)
Write-Host "Micro results cell updated. Paragraph count:" $cell.Range.Paragraphs.Count
